$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 13.13709999999999
$ws.Range("E6").Value = 11.5955
$ws.Range("D7").Value = -7.518199999999992
$ws.Range("E7").Value = 12.7404
$ws.Range("A8").Value = -21.19700000000001
$ws.Range("E8").Value = 13.0751
$ws.Range("E9").Value = 9.195899999999988
$ws.Range("A10").Value = -20.46869999999997
$ws.Range("E10").Value = 11.1834
$ws.Range("A12").Value = -22.57200000000002
$ws.Range("E12").Value = 12.66269999999999
$ws.Range("C13").Value = -12.97589999999999
$ws.Range("A18").Value = -22.42220000000002
$ws.Range("D20").Value = -8.315499999999997
$ws.Range("A25").Value = -22.31520000000003
